$d = $word.ActiveDocument

# Locate the "Activities" bullet list (numId 11) anchored by its first item's
# distinctive original text, and the "Notes" sub-bullets (numId 12) that get
# removed entirely. We find the paragraph indices by scanning for the known
# text so the script is resilient to any surrounding content.
$count = $d.Paragraphs.Count

$idxCredits = -1
$idxCreateYourVideo = -1
$idxMay2018 = -1
$idxInFuture = -1

for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($idxCredits -eq -1 -and $t -like "Tap to insert local credits*") {
        $idxCredits = $i
    }
    if ($idxCreateYourVideo -eq -1 -and $t -like "Create your new video.*") {
        $idxCreateYourVideo = $i
    }
    if ($idxMay2018 -eq -1 -and $t -like "May 2018:*") {
        $idxMay2018 = $i
    }
    if ($idxInFuture -eq -1 -and $t -like "In the future, some file naming features will be automated.*") {
        $idxInFuture = $i
    }
}

# The five "Activities" bullet paragraphs that shift their text up by one slot.
$p1 = $idxCredits
$p2 = $p1 + 1
$p3 = $p1 + 2
$p4 = $p1 + 3
$p5 = $p1 + 4
$p6 = $idxCreateYourVideo

# Capture the original texts (without trailing paragraph marks) before any edits.
$text2 = $d.Paragraphs($p2).Range.Text
$text3 = $d.Paragraphs($p3).Range.Text
$text4 = $d.Paragraphs($p4).Range.Text
$text5 = $d.Paragraphs($p5).Range.Text
$text6 = $d.Paragraphs($p6).Range.Text

function Strip-ParaMark($s) {
    return $s.TrimEnd([char]13, [char]7)
}

$text2 = Strip-ParaMark $text2
$text3 = Strip-ParaMark $text3
$text4 = Strip-ParaMark $text4
$text5 = Strip-ParaMark $text5
$text6 = Strip-ParaMark $text6

# Shift each bullet's text into the previous bullet, recoloring the first one
# from red to black to indicate it is now a finished (reviewed) instruction.
$d.Paragraphs($p1).Range.Text = $text2
$d.Paragraphs($p1).Range.Font.Color = 0
$d.Paragraphs($p2).Range.Text = $text3
$d.Paragraphs($p3).Range.Text = $text4
$d.Paragraphs($p4).Range.Text = $text5
$d.Paragraphs($p5).Range.Text = $text6

# Remove the now-duplicated trailing bullet paragraph entirely.
$d.Paragraphs($p6).Range.Delete()

# Remove the four red "Notes" sub-bullets (May 2018 ... In the future ...)
# as a block, now that the document shifted by one fewer paragraph.
$idxMay2018b = -1
$idxInFutureb = -1
$count2 = $d.Paragraphs.Count
for ($i = 1; $i -le $count2; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($idxMay2018b -eq -1 -and $t -like "May 2018:*") {
        $idxMay2018b = $i
    }
    if ($idxInFutureb -eq -1 -and $t -like "In the future, some file naming features will be automated.*") {
        $idxInFutureb = $i
    }
}

$startRange = $d.Paragraphs($idxMay2018b).Range.Start
$endRange = $d.Paragraphs($idxInFutureb).Range.End
$blockRange = $d.Range($startRange, $endRange)
$blockRange.Delete()
